$d = $word.ActiveDocument

function Replace-All($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

function Replace-One($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 1) | Out-Null
}

# --- Language switcher line at top of document ---
# "English" appears twice (link text + label further down); both become "Inglês"
Replace-All "English" "Inglês"
Replace-All " / Portuguese / French / Thai / Vietnamese / Spanish" " / Português / Francês / Tailandês / Vietnamita / Espanhol"

# --- Brief / target audience table ---
Replace-All "Brief" "Resumo"
Replace-All "An email to partners in the target country who haven’t RSVPed to remind them to send the RSVP. It will be sent via customer.io" "Um e-mail para parceiros no país-alvo que ainda não enviaram o RSVP para os lembrar de enviar a resposta. Será enviado através do customer.io"
Replace-All "Target audience" "Público-alvo"
Replace-All "Invited partners who haven’t RSVPed yet" "Parceiros convidados que ainda não enviaram o RSVP"

# --- Subject line ---
Replace-All "Subject line" "Linha de assunto"
Replace-All ": Reminder: RSVP for " ": Reserve já o seu lugar para o "

# --- [EVENT NAME] placeholder (5 uniform occurrences, incl. one with trailing space) ---
Replace-All "[EVENT NAME]" "[NOME DO EVENTO]"

# --- Headline ---
Replace-All "Don’t delay! Book your spot today!" "Não adie! Reserve já o seu lugar!"

# --- Greeting ---
Replace-All "Hi " "Olá "
Replace-All "[PARTNER NAME]" "[NOME DO PARCEIRO]"

# --- Event excitement paragraphs (one-day and multi-day variants) ---
Replace-All "We hope you’re as excited as us for the " "Esperamos que esteja tão entusiasmado quanto nós com o "
Replace-All ", happening on " ", que irá decorrer no dia "
Replace-All ", happening from " ", a decorrer de "
Replace-One " to " " a "

# --- [DD Mmm YYYY] / DD Mmm YYYY (4 uniform occurrences, bracketed and bare) ---
Replace-All "DD Mmm YYYY" "DD Mmm AAAA"

# --- Confirm attendance paragraph ---
Replace-All "Confirm your attendance for this highly-anticipated event by [" "Confirme a sua presença neste evento tão aguardado até ["
Replace-All "] as spots are limited and on a first-come, first-served basis." "], uma vez que as vagas são limitadas e por ordem de chegada."

# --- RSVP button ---
Replace-All "RVSP now" "Preencher formulário"

# --- Contact paragraphs ---
Replace-All "If you have any questions, please contact us via " "Para mais informações, contacte-nos através do "
Replace-All " or " " ou "
Replace-All "If you have any questions, please contact your country manager, " "Para mais questões, pode também contactar o seus gestor de parcerias "
Replace-All ", at " ", em "

# --- Look forward to seeing you (two variants) ---
Replace-All "We look forward to seeing you at " "Esperamos vê-lo em breve, no evento "
Replace-All "We look forward to seeing you at [EVENT NAME]! " "Esperamos vê-lo em breve, no evento [NOME DO EVENTO]! "

# --- Numbered list items (centered footer block) ---
Replace-All "If you have any questions, please contact your country manager:" "Se tiver alguma dúvida, contacte o gestor do seu país:"
Replace-All " [NAME] |  [EMAIL ADDRESS] | [WHATSAPP NO] (WhatsApp). " " [NAME] | [ENDEREÇO DE EMAIL] | [NO DO WHATSAPP] (WhatsApp). "
Replace-All "If you have any questions, please contact us via:" "Se tiver alguma dúvida, contacte-nos através de:"

# --- Comments (the only editable surface for comment text rewrites the comment's
#     paragraph/run formatting to a default style in this runtime, but keeps the
#     textual content correct) ---
$c0 = $d.Comments.Item(1)
$c0.Content = "escolha a primeira opção se for um evento de um dia`r`r" + `
              "escolha a segunda opção se for um evento de vários dias"

$d.Comments.Item(2).Content = "Escolha um deles"
$d.Comments.Item(3).Content = "Escolha um deles"
$d.Comments.Item(4).Content = "Escolha um deles"
